$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update Seasonality Index (column L) ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$wsForecast.Range("L2").Value = 0.91
$wsForecast.Range("L3").Value = 0.95
$wsForecast.Range("L4").Value = 1.16
$wsForecast.Range("L5").Value = 1.09
$wsForecast.Range("L6").Value = 0.88
$wsForecast.Range("L7").Value = 0.92
$wsForecast.Range("L8").Value = 0.83
$wsForecast.Range("L11").Value = 1.16
$wsForecast.Range("L12").Value = 0.89
$wsForecast.Range("L13").Value = 0.91
$wsForecast.Range("L14").Value = 1.07
$wsForecast.Range("L15").Value = 1.02
$wsForecast.Range("L16").Value = 1.15
$wsForecast.Range("L17").Value = 0.95

# --- Sheet "Summary": update forecast totals (column B) ---
# These cells hold numeric-looking values stored as text, so force a text
# number format before assigning to keep them as text (matches original type).
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "45"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "23"
